$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"   = -8.590900000000001
    "E3"   = 16.1486
    "B12"  = 4.9978
    "D14"  = -7.869399999999996
    "D26"  = -8.784500000000012
    "E30"  = 15.76270000000001
    "D31"  = -9.020400000000002
    "B32"  = 6.502099999999998
    "D35"  = -8.329899999999999
    "B36"  = 9.2759
    "D37"  = -7.636199999999998
    "B38"  = 4.852899999999998
    "E44"  = 16.51669999999999
    "D45"  = -7.656499999999998
    "B46"  = 6.097800000000001
    "B54"  = 4.488899999999997
    "B55"  = 5.530999999999998
    "D57"  = -8.398999999999994
    "E58"  = 16.33120000000001
    "B67"  = 6.260099999999997
    "B69"  = 5.416099999999997
    "B72"  = 5.169800000000004
    "E84"  = 16.54649999999999
    "E89"  = 17.34830000000002
    "B91"  = 5.833900000000001
    "E91"  = 18.08000000000002
    "E92"  = 18.09600000000002
    "B99"  = 4.600499999999997
    "D100" = -8.769499999999994
    "D102" = -7.795400000000001
    "E102" = 16.5073
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
